$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------------
# New Google-Forms quiz submissions arrived (rows 518-535); row 517 stops being
# the last row of the table, so its "final row" styling moves down to the new
# last row (535) and 517 takes on the regular interior-row style.
# ----------------------------------------------------------------------------

# 1) Snapshot row 517s current ("last row") styling into new row 535 first,
#    before row 517 itself gets restyled.
$ws.Range("A517:N517").Copy($ws.Range("A535:N535"))

# 2) Row 517 becomes a normal interior row (style-wise).
$ws.Range("A515:N515").Copy($ws.Range("A517:N517"))

# 3) Stamp interior alternating-row styles onto the newly appended rows 518-534.
$ws.Range("A514:N514").Copy($ws.Range("A518:N518"))
$ws.Range("A515:N515").Copy($ws.Range("A519:N519"))
$ws.Range("A516:N516").Copy($ws.Range("A520:N520"))
$ws.Range("A515:N515").Copy($ws.Range("A521:N521"))
$ws.Range("A516:N516").Copy($ws.Range("A522:N522"))
$ws.Range("A503:N503").Copy($ws.Range("A523:N523"))
$ws.Range("A516:N516").Copy($ws.Range("A524:N524"))
$ws.Range("A503:N503").Copy($ws.Range("A525:N525"))
$ws.Range("A516:N516").Copy($ws.Range("A526:N526"))
$ws.Range("A515:N515").Copy($ws.Range("A527:N527"))
$ws.Range("A516:N516").Copy($ws.Range("A528:N528"))
$ws.Range("A503:N503").Copy($ws.Range("A529:N529"))
$ws.Range("A514:N514").Copy($ws.Range("A530:N530"))
$ws.Range("A503:N503").Copy($ws.Range("A531:N531"))
$ws.Range("A514:N514").Copy($ws.Range("A532:N532"))
$ws.Range("A515:N515").Copy($ws.Range("A533:N533"))
$ws.Range("A516:N516").Copy($ws.Range("A534:N534"))

# 4) Each Copy() leaves a phantom blank M/N cell behind (the template row used
#    the opposite "extra" column) -- drop those so only one of M/N exists,
#    matching this sheets convention of never emitting both.
$ws.Range("M518").ClearContents()
$ws.Range("N519").ClearContents()
$ws.Range("N520").ClearContents()
$ws.Range("N521").ClearContents()
$ws.Range("N522").ClearContents()
$ws.Range("M523").ClearContents()
$ws.Range("N524").ClearContents()
$ws.Range("M525").ClearContents()
$ws.Range("N526").ClearContents()
$ws.Range("N527").ClearContents()
$ws.Range("N528").ClearContents()
$ws.Range("M529").ClearContents()
$ws.Range("M530").ClearContents()
$ws.Range("M531").ClearContents()
$ws.Range("M532").ClearContents()
$ws.Range("N533").ClearContents()
$ws.Range("N534").ClearContents()
$ws.Range("N535").ClearContents()

# 5) Write the actual submitted values for rows 517 (unchanged) through 535 (new).
# Row 517
$ws.Cells.Item(517, 1).Value = 45569.873526111114
$ws.Cells.Item(517, 2).Value = 'kty030122@gmail.com'
$ws.Cells.Item(517, 3).Value = '체육학과'
$ws.Cells.Item(517, 4).Value = 20214113.0
$ws.Cells.Item(517, 5).Value = '김태연'
$ws.Cells.Item(517, 6).Value = '‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다.'
$ws.Cells.Item(517, 7).Value = 0.5
$ws.Cells.Item(517, 8).Value = '6:4'
$ws.Cells.Item(517, 9).Value = '15분의 1'
$ws.Cells.Item(517, 10).Value = '44만호, 153만명'
$ws.Cells.Item(517, 11).Value = '경상'
$ws.Cells.Item(517, 12).Value = 'Red'
$ws.Cells.Item(517, 13).Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# Row 518
$ws.Cells.Item(518, 1).Value = 45569.92259511574
$ws.Cells.Item(518, 2).Value = 'goeunsue@naver.com'
$ws.Cells.Item(518, 3).Value = '경영대학'
$ws.Cells.Item(518, 4).Value = 20242907.0
$ws.Cells.Item(518, 5).Value = '고은수'
$ws.Cells.Item(518, 6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(518, 7).Value = 0.1
$ws.Cells.Item(518, 8).Value = '6:4'
$ws.Cells.Item(518, 9).Value = '20분의 1'
$ws.Cells.Item(518, 10).Value = '20만호, 69만명'
$ws.Cells.Item(518, 11).Value = '충청'
$ws.Cells.Item(518, 12).Value = 'Black'
$ws.Cells.Item(518, 14).Value = '모름/무응답'

# Row 519
$ws.Cells.Item(519, 1).Value = 45569.95132994213
$ws.Cells.Item(519, 2).Value = 'taewon16@naver.com'
$ws.Cells.Item(519, 3).Value = '빅데이터학과'
$ws.Cells.Item(519, 4).Value = 20195158.0
$ws.Cells.Item(519, 5).Value = '류태원'
$ws.Cells.Item(519, 6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(519, 7).Value = 0.1
$ws.Cells.Item(519, 8).Value = '6:4'
$ws.Cells.Item(519, 9).Value = '15분의 1'
$ws.Cells.Item(519, 10).Value = '20만호, 69만명'
$ws.Cells.Item(519, 11).Value = '충청'
$ws.Cells.Item(519, 12).Value = 'Red'
$ws.Cells.Item(519, 13).Value = '반대한다.'

# Row 520
$ws.Cells.Item(520, 1).Value = 45569.96054436342
$ws.Cells.Item(520, 2).Value = 'jb9517asd@naver.com'
$ws.Cells.Item(520, 3).Value = '소프트웨어학부'
$ws.Cells.Item(520, 4).Value = 20245109.0
$ws.Cells.Item(520, 5).Value = '곽우주'
$ws.Cells.Item(520, 6).Value = '과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다.'
$ws.Cells.Item(520, 7).Value = 0.7
$ws.Cells.Item(520, 8).Value = '7:3'
$ws.Cells.Item(520, 9).Value = '10분의 1'
$ws.Cells.Item(520, 10).Value = '130만호, 5백만명'
$ws.Cells.Item(520, 11).Value = '경기'
$ws.Cells.Item(520, 12).Value = 'Red'
$ws.Cells.Item(520, 13).Value = '반대한다.'

# Row 521
$ws.Cells.Item(521, 1).Value = 45569.96751702546
$ws.Cells.Item(521, 2).Value = '1202kge@naver.com'
$ws.Cells.Item(521, 3).Value = '사회학과'
$ws.Cells.Item(521, 4).Value = 20242205.0
$ws.Cells.Item(521, 5).Value = '김가은'
$ws.Cells.Item(521, 6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(521, 7).Value = 0.7
$ws.Cells.Item(521, 8).Value = '4:6'
$ws.Cells.Item(521, 9).Value = '10분의 1'
$ws.Cells.Item(521, 10).Value = '44만호, 153만명'
$ws.Cells.Item(521, 11).Value = '전라'
$ws.Cells.Item(521, 12).Value = 'Red'
$ws.Cells.Item(521, 13).Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# Row 522
$ws.Cells.Item(522, 1).Value = 45569.97317773148
$ws.Cells.Item(522, 2).Value = 'sowon051125@naver.com'
$ws.Cells.Item(522, 3).Value = '데이터사이언스학부'
$ws.Cells.Item(522, 4).Value = 20243238.0
$ws.Cells.Item(522, 5).Value = '이소원'
$ws.Cells.Item(522, 6).Value = '‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다.'
$ws.Cells.Item(522, 7).Value = 0.3
$ws.Cells.Item(522, 8).Value = '5:5'
$ws.Cells.Item(522, 9).Value = '20분의 1'
$ws.Cells.Item(522, 10).Value = '15만호,  32만명'
$ws.Cells.Item(522, 11).Value = '경상'
$ws.Cells.Item(522, 12).Value = 'Red'
$ws.Cells.Item(522, 13).Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# Row 523
$ws.Cells.Item(523, 1).Value = 45569.982878310184
$ws.Cells.Item(523, 2).Value = 'dlxotjq27@gmail.com'
$ws.Cells.Item(523, 3).Value = '경영학과'
$ws.Cells.Item(523, 4).Value = 20213023.0
$ws.Cells.Item(523, 5).Value = '이태섭'
$ws.Cells.Item(523, 6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(523, 7).Value = 0.1
$ws.Cells.Item(523, 8).Value = '6:4'
$ws.Cells.Item(523, 9).Value = '20분의 1'
$ws.Cells.Item(523, 10).Value = '20만호, 69만명'
$ws.Cells.Item(523, 11).Value = '충청'
$ws.Cells.Item(523, 12).Value = 'Black'
$ws.Cells.Item(523, 14).Value = '노동자가 과도한 연장근로를 받을 수 있어 반대한다.'

# Row 524
$ws.Cells.Item(524, 1).Value = 45569.984793599535
$ws.Cells.Item(524, 2).Value = 'jytoto33@naver.com'
$ws.Cells.Item(524, 3).Value = '언어청각학부'
$ws.Cells.Item(524, 4).Value = 20243912.0
$ws.Cells.Item(524, 5).Value = '김지윤'
$ws.Cells.Item(524, 6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(524, 7).Value = 0.3
$ws.Cells.Item(524, 8).Value = '3:7'
$ws.Cells.Item(524, 9).Value = '30분의 1'
$ws.Cells.Item(524, 10).Value = '44만호, 153만명'
$ws.Cells.Item(524, 11).Value = '평안'
$ws.Cells.Item(524, 12).Value = 'Red'
$ws.Cells.Item(524, 13).Value = '반대한다.'

# Row 525
$ws.Cells.Item(525, 1).Value = 45570.011931076384
$ws.Cells.Item(525, 2).Value = 'andy041001@naver.com'
$ws.Cells.Item(525, 3).Value = '러시아학과'
$ws.Cells.Item(525, 4).Value = 20231720.0
$ws.Cells.Item(525, 5).Value = '이형범'
$ws.Cells.Item(525, 6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(525, 7).Value = 0.9
$ws.Cells.Item(525, 8).Value = '4:6'
$ws.Cells.Item(525, 9).Value = '10분의 1'
$ws.Cells.Item(525, 10).Value = '44만호, 153만명'
$ws.Cells.Item(525, 11).Value = '전라'
$ws.Cells.Item(525, 12).Value = 'Black'
$ws.Cells.Item(525, 14).Value = '노동자가 과도한 연장근로를 받을 수 있어 반대한다.'

# Row 526
$ws.Cells.Item(526, 1).Value = 45570.02593584491
$ws.Cells.Item(526, 2).Value = 'hyerim0v0@gmail.com'
$ws.Cells.Item(526, 3).Value = '일본학과'
$ws.Cells.Item(526, 4).Value = 20231630.0
$ws.Cells.Item(526, 5).Value = '전혜림'
$ws.Cells.Item(526, 6).Value = '실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다.'
$ws.Cells.Item(526, 7).Value = 0.3
$ws.Cells.Item(526, 8).Value = '6:4'
$ws.Cells.Item(526, 9).Value = '15분의 1'
$ws.Cells.Item(526, 10).Value = '20만호, 69만명'
$ws.Cells.Item(526, 11).Value = '경상'
$ws.Cells.Item(526, 12).Value = 'Red'
$ws.Cells.Item(526, 13).Value = '반대한다.'

# Row 527
$ws.Cells.Item(527, 1).Value = 45570.040854652776
$ws.Cells.Item(527, 2).Value = 'kby5432@naver.com'
$ws.Cells.Item(527, 3).Value = '법학과'
$ws.Cells.Item(527, 4).Value = 20192737.0
$ws.Cells.Item(527, 5).Value = '윤경빈'
$ws.Cells.Item(527, 6).Value = '‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다.'
$ws.Cells.Item(527, 7).Value = 0.1
$ws.Cells.Item(527, 8).Value = '7:3'
$ws.Cells.Item(527, 9).Value = '15분의 1'
$ws.Cells.Item(527, 10).Value = '44만호, 153만명'
$ws.Cells.Item(527, 11).Value = '평안'
$ws.Cells.Item(527, 12).Value = 'Red'
$ws.Cells.Item(527, 13).Value = '반대한다.'

# Row 528
$ws.Cells.Item(528, 1).Value = 45570.04539924768
$ws.Cells.Item(528, 2).Value = 'jign1106@naver.com'
$ws.Cells.Item(528, 3).Value = '간호학과'
$ws.Cells.Item(528, 4).Value = 20246289.0
$ws.Cells.Item(528, 5).Value = '지은총'
$ws.Cells.Item(528, 6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(528, 7).Value = 0.1
$ws.Cells.Item(528, 8).Value = '6:4'
$ws.Cells.Item(528, 9).Value = '20분의 1'
$ws.Cells.Item(528, 10).Value = '20만호, 69만명'
$ws.Cells.Item(528, 11).Value = '충청'
$ws.Cells.Item(528, 12).Value = 'Red'
$ws.Cells.Item(528, 13).Value = '반대한다.'

# Row 529
$ws.Cells.Item(529, 1).Value = 45570.080389097224
$ws.Cells.Item(529, 2).Value = 'kt433@naver.com'
$ws.Cells.Item(529, 3).Value = '사회복지학과'
$ws.Cells.Item(529, 4).Value = 20222361.0
$ws.Cells.Item(529, 5).Value = '주혜린'
$ws.Cells.Item(529, 6).Value = '‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다.'
$ws.Cells.Item(529, 7).Value = 0.1
$ws.Cells.Item(529, 8).Value = '7:3'
$ws.Cells.Item(529, 9).Value = '10분의 1'
$ws.Cells.Item(529, 10).Value = '20만호, 69만명'
$ws.Cells.Item(529, 11).Value = '충청'
$ws.Cells.Item(529, 12).Value = 'Black'
$ws.Cells.Item(529, 14).Value = '노동자가 과도한 연장근로를 받을 수 있어 반대한다.'

# Row 530
$ws.Cells.Item(530, 1).Value = 45570.13000396991
$ws.Cells.Item(530, 2).Value = 'misunhong0707@gmail.com'
$ws.Cells.Item(530, 3).Value = '융합과학수사학과'
$ws.Cells.Item(530, 4).Value = 20246940.0
$ws.Cells.Item(530, 5).Value = '홍미선'
$ws.Cells.Item(530, 6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(530, 7).Value = 0.1
$ws.Cells.Item(530, 8).Value = '6:4'
$ws.Cells.Item(530, 9).Value = '20분의 1'
$ws.Cells.Item(530, 10).Value = '20만호, 69만명'
$ws.Cells.Item(530, 11).Value = '경상'
$ws.Cells.Item(530, 12).Value = 'Black'
$ws.Cells.Item(530, 14).Value = '노동자가 과도한 연장근로를 받을 수 있어 반대한다.'

# Row 531
$ws.Cells.Item(531, 1).Value = 45570.133293969906
$ws.Cells.Item(531, 2).Value = 'ziva0726@naver.com'
$ws.Cells.Item(531, 3).Value = '심리학과'
$ws.Cells.Item(531, 4).Value = 20212104.0
$ws.Cells.Item(531, 5).Value = '김소현'
$ws.Cells.Item(531, 6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(531, 7).Value = 0.3
$ws.Cells.Item(531, 8).Value = '6:4'
$ws.Cells.Item(531, 9).Value = '30분의 1'
$ws.Cells.Item(531, 10).Value = '20만호, 69만명'
$ws.Cells.Item(531, 11).Value = '전라'
$ws.Cells.Item(531, 12).Value = 'Black'
$ws.Cells.Item(531, 14).Value = '노동자가 과도한 연장근로를 받을 수 있어 반대한다.'

# Row 532
$ws.Cells.Item(532, 1).Value = 45570.14750028936
$ws.Cells.Item(532, 2).Value = 'sujdiamond@gmail.com'
$ws.Cells.Item(532, 3).Value = '바이오메디컬학과'
$ws.Cells.Item(532, 4).Value = 20243627.0
$ws.Cells.Item(532, 5).Value = '심유진'
$ws.Cells.Item(532, 6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(532, 7).Value = 0.1
$ws.Cells.Item(532, 8).Value = '6:4'
$ws.Cells.Item(532, 9).Value = '20분의 1'
$ws.Cells.Item(532, 10).Value = '20만호, 69만명'
$ws.Cells.Item(532, 11).Value = '충청'
$ws.Cells.Item(532, 12).Value = 'Black'
$ws.Cells.Item(532, 14).Value = '찬성한다.'

# Row 533
$ws.Cells.Item(533, 1).Value = 45570.168118692134
$ws.Cells.Item(533, 2).Value = '20182346@hallym.ac.kr'
$ws.Cells.Item(533, 3).Value = '사회복지학부'
$ws.Cells.Item(533, 4).Value = 20182346.0
$ws.Cells.Item(533, 5).Value = '이용재'
$ws.Cells.Item(533, 6).Value = '‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다.'
$ws.Cells.Item(533, 7).Value = 0.1
$ws.Cells.Item(533, 8).Value = '7:3'
$ws.Cells.Item(533, 9).Value = '20분의 1'
$ws.Cells.Item(533, 10).Value = '130만호, 5백만명'
$ws.Cells.Item(533, 11).Value = '전라'
$ws.Cells.Item(533, 12).Value = 'Red'
$ws.Cells.Item(533, 13).Value = '반대한다.'

# Row 534
$ws.Cells.Item(534, 1).Value = 45570.17452395833
$ws.Cells.Item(534, 2).Value = 'dncks5343@naver.com'
$ws.Cells.Item(534, 3).Value = '언어청각학부'
$ws.Cells.Item(534, 4).Value = 20243973.0
$ws.Cells.Item(534, 5).Value = '장우찬'
$ws.Cells.Item(534, 6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(534, 7).Value = 0.1
$ws.Cells.Item(534, 8).Value = '6:4'
$ws.Cells.Item(534, 9).Value = '20분의 1'
$ws.Cells.Item(534, 10).Value = '20만호, 69만명'
$ws.Cells.Item(534, 11).Value = '충청'
$ws.Cells.Item(534, 12).Value = 'Red'
$ws.Cells.Item(534, 13).Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# Row 535
$ws.Cells.Item(535, 1).Value = 45570.24906371528
$ws.Cells.Item(535, 2).Value = 'leedongyoung797@gmail.com'
$ws.Cells.Item(535, 3).Value = '언어청각학부'
$ws.Cells.Item(535, 4).Value = 20243934.0
$ws.Cells.Item(535, 5).Value = '이동영'
$ws.Cells.Item(535, 6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(535, 7).Value = 0.3
$ws.Cells.Item(535, 8).Value = '6:4'
$ws.Cells.Item(535, 9).Value = '10분의 1'
$ws.Cells.Item(535, 10).Value = '20만호, 69만명'
$ws.Cells.Item(535, 11).Value = '전라'
$ws.Cells.Item(535, 12).Value = 'Red'
$ws.Cells.Item(535, 13).Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# 6) The response table now spans through the new last row.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:N535"))
